$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.759.31'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '2.525.49'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.27'
$ws.Range('E5').Value = '  +1.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.50'
$ws.Range('E6').Value = '  +1.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.573'
$ws.Range('E7').Value = '  -0.83%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.535'
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.67'
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.48'
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('E13').Value = '  -4.13%  '
$ws.Range('D14').Value = '2.915.70'
$ws.Range('E14').Value = '  -0.48%  '
$ws.Range('D15').Value = '2.551.85'
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.01'
$ws.Range('E16').Value = '  -4.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.847'
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('D18').Value = '42.829.02'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.86'
$ws.Range('E19').Value = '  +2.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.49'
$ws.Range('E20').Value = '  -4.69%  '
$ws.Range('D21').Value = '0.0₃0962'
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.59'
$ws.Range('E22').Value = '  -2.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '251.80'
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.96'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.06'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.28'
$ws.Range('E26').Value = '  -5.09%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.42'
$ws.Range('E28').Value = '  +2.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.99'
$ws.Range('E29').Value = '  +2.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.41'
$ws.Range('E30').Value = '  +3.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.86'
$ws.Range('E31').Value = '  -2.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '156.40'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('B33').Value = 'Celestia'
$ws.Range('C33').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.46'
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.14'
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('E35').Value = '  -1.46%  '
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0788'
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('E38').Value = '  +0.99%  '
$ws.Range('E39').Value = '  +10.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.118'
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '21.64'
$ws.Range('E41').Value = '  -13.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0304'
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.79'
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('E45').Value = '  -3.85%  '
$ws.Range('D46').Value = '1.994.64'
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.12'
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '83.92'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.51'
$ws.Range('E49').Value = '  +3.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.84'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('D51').Value = '2.773.89'
$ws.Range('E51').Value = '  -0.41%  '
